$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prompt = @"
You are an assistant for teaching and learning English. Your task is to classify the intent of the user's utterance based on the intent list provided in the question
<task>
Step 1: Read the assistant's question and user's utterance.
Step 2: Based on the information describing the customer's intent, determine which intent category the answer belongs to from the list below. If it doesn't match any intent, classify it as a "fallback" intent.
Step 3: Return the intent name.
</task>
<tag>
## Descripble intent list
[
        {
                "intent_name": "affirm_confirm",
                "description": "Khách hàng nói xác nhận đúng khách hàng"
        },
        {
                "intent_name": "absent",
                "description": "Khách hàng nói đi vắng hoặc nghe hộ"
        },
        {
                "intent_name": "deny_confirm",
                "description": "Khách hàng xác nhận không phải khách hàng hoặc nhầm máy"
        },
        {
                "intent_name": "busy",
                "description": "khách hàng báo bận hoặc đang họp"
        },
        {
                "intent_name": "cant_hear",
                "description": "Khách hàng nói không nghe rõ"
        },
        {
                "intent_name": "fallback",
                "description": "Khách hàng nói những câu nói ngoài phạm vi không liên quan tới câu hỏi, nói tục hoặc chửi bậy"
        }
]
</tag>
<ouput>
The result should return only one intent that best matches the customer's response.
The returned intent must belong to one of the intent lists mentioned above.
Only the intent name should be generated, no other characters are allowed.
</ouput>
"@

$apiKeyError = "Request failed after 2 retries: Error code: 401 - {'error': {'message': 'Invalid API Key', 'type': 'invalid_request_error', 'code': 'invalid_api_key'}}"

# --- Update existing row 2 (response_time changed) ---
$ws.Cells.Item(2, 6).Value = 1.109622955322266

# --- Update existing row 3 (now the groq "Hello." failed-auth row) ---
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "groq:llama-3.3-70b-versatile"
$ws.Cells.Item(3, 4).Value = "Hello."
$ws.Cells.Item(3, 5).Value = $apiKeyError
$ws.Cells.Item(3, 6).Value = -1

# --- Row 4: openai:gpt-4o-mini / "I eat breakfast." success ---
$ws.Cells.Item(4, 1).Value = ""
$ws.Cells.Item(4, 2).Value = "openai:gpt-4o-mini"
$ws.Cells.Item(4, 3).Value = $prompt
$ws.Cells.Item(4, 4).Value = "I eat breakfast."
$ws.Cells.Item(4, 5).Value = "fallback"
$ws.Cells.Item(4, 6).Value = 0.6955010890960693

# --- Row 5: groq:llama-3.3-70b-versatile / "I eat breakfast." failed auth ---
$ws.Cells.Item(5, 1).Value = ""
$ws.Cells.Item(5, 2).Value = "groq:llama-3.3-70b-versatile"
$ws.Cells.Item(5, 3).Value = $prompt
$ws.Cells.Item(5, 4).Value = "I eat breakfast."
$ws.Cells.Item(5, 5).Value = $apiKeyError
$ws.Cells.Item(5, 6).Value = -1

# --- Row 6: openai:gpt-4o-mini / "I eat breakfast." success ---
$ws.Cells.Item(6, 1).Value = ""
$ws.Cells.Item(6, 2).Value = "openai:gpt-4o-mini"
$ws.Cells.Item(6, 3).Value = $prompt
$ws.Cells.Item(6, 4).Value = "I eat breakfast."
$ws.Cells.Item(6, 5).Value = "fallback"
$ws.Cells.Item(6, 6).Value = 0.7523543834686279

# --- Row 7: groq:llama-3.3-70b-versatile / "I eat breakfast." failed auth ---
$ws.Cells.Item(7, 1).Value = ""
$ws.Cells.Item(7, 2).Value = "groq:llama-3.3-70b-versatile"
$ws.Cells.Item(7, 3).Value = $prompt
$ws.Cells.Item(7, 4).Value = "I eat breakfast."
$ws.Cells.Item(7, 5).Value = $apiKeyError
$ws.Cells.Item(7, 6).Value = -1
